{"js": "// Remove the \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph\n// together with the blank paragraph right before it, the blank paragraph\n// right after it, and the (empty) page-break paragraph that follows that,\n// exactly as removed by the source diff.\n\nconst body = context.document.body;\n\n// Locate the distinctive paragraph by its text.\nconst results = body.search(\"Ver no Jupiter Salvar em pdf Salvar em docx\", {\n  matchCase: true\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Paragraph \"Ver no Jupiter Salvar em pdf Salvar em docx\" not found.');\n}\n\nconst target = results.items[0].paragraphs.getFirst();\nconst before = target.getPrevious();   // blank paragraph just above the target\nconst after1 = target.getNext();       // blank paragraph just below the target\nconst after2 = after1.getNext();       // empty pageBreakBefore paragraph after that\nawait context.sync();\n\n// Delete from bottom to top so earlier references stay valid.\nafter2.delete();\nafter1.delete();\ntarget.delete();\nbefore.delete();\n\nawait context.sync();\n", "ps1": "# Remove the \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph\n# together with the blank paragraph right before it, the blank paragraph\n# right after it, and the (empty) page-break paragraph that follows that,\n# exactly as removed by the source diff.\n\n$d = $word.ActiveDocument\n\n# Confirm the target text is present in the document (per the expected edit).\n$found = $d.Content.Find.Execute(\"Ver no Jupiter Salvar em pdf Salvar em docx\")\nif (-not $found) {\n    throw \"Target paragraph text not found in document.\"\n}\n\n# Locate the exact paragraph index (Find's match range does not reliably\n# expose the owning paragraph in this host, so resolve it via the\n# Paragraphs collection instead).\n$target = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd()\n    if ($text -eq \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n        $target = $i\n        break\n    }\n}\nif ($target -eq -1) {\n    throw \"Target paragraph not found.\"\n}\n\n# Delete the 4 paragraphs: the one before, the target, and the two after it.\n# Delete from the bottom up so earlier indices stay valid.\n$d.Paragraphs.Item($target + 2).Range.Delete()\n$d.Paragraphs.Item($target + 1).Range.Delete()\n$d.Paragraphs.Item($target).Range.Delete()\n$d.Paragraphs.Item($target - 1).Range.Delete()\n"}
